$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.549.84'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.846.68'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '264.07'
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '0.5216'
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('D8').Value = '0.3225'
$ws.Range('E8').Value = '  -0.74%  '
$ws.Range('D9').Value = '0.06803'
$ws.Range('E9').Value = '  +0.42%  '
$ws.Range('D10').Value = '18.79'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').Value = '0.7790'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('D12').Value = '0.07764'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '1.823.50'
$ws.Range('E13').Value = '  -1.73%  '
$ws.Range('D14').Value = '88.50'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('D15').Value = '5.026'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Value = '13.98'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '0.000007966'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D20').Value = '26.596.88'
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('D21').Value = '4.628'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').Value = '9.462'
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('D23').Value = '6.006'
$ws.Range('E23').Value = '  +1.39%  '
$ws.Range('D24').Value = '143.18'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('D25').Value = '2.175'
$ws.Range('E25').Value = '  -7.51%  '
$ws.Range('D26').Value = '1.681'
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('D28').Value = '111.77'
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('D29').Value = '4.184'
$ws.Range('E29').Value = '  -0.79%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.08740'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.115'
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('D32').Value = '0.04840'
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').Value = '0.7208'
$ws.Range('E33').Value = '  +4.75%  '
$ws.Range('D34').Value = '1.129'
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('D35').Value = '2.862'
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').Value = '0.01792'
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('D39').Value = '0.4862'
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '0.8965'
$ws.Range('E40').Value = '  -0.51%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').Value = '110.94'
$ws.Range('E41').Value = '  -1.79%  '
$ws.Range('D42').Value = '6.029'
$ws.Range('E42').Value = '  -1.88%  '
$ws.Range('D43').Value = '0.9999'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').Value = '7.636'
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('D45').Value = '0.4205'
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('D46').Value = '0.05890'
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '9.032'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('D49').Value = '35.00'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').Value = '0.8877'
$ws.Range('E50').Value = '  +3.54%  '
$ws.Range('D51').Value = '59.85'
$ws.Range('E51').Value = '  +0.85%  '
